# Insert a new data row before the current row 198, shifting the existing
# row 198 (and everything after it, through the old row 292) down by one.
# Excel's native row-insert semantics handle the shift of all columns
# (A:R) automatically; we then populate the newly-inserted row 198 with
# the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(198).Insert()

$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = 44917
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = 100112039
$ws.Cells.Item(198, 7).Value = "Ciboulette"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 120
$ws.Cells.Item(198, 11).Value = 6000
$ws.Cells.Item(198, 12).Value = 6000
$ws.Cells.Item(198, 13).Value = 6000
$ws.Cells.Item(198, 14).Value = "`$/docena de atados"
$ws.Cells.Item(198, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(198, 16).Value = 2000
$ws.Cells.Item(198, 17).Value = 3
$ws.Cells.Item(198, 18).Value = "Hortaliza"
